$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.249.53'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '''3.942.33'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''498.04'
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("D6").Value = '''148.19'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '''0.625'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").Value = '''0.998'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.735'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("D10").Value = '''0.177'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = '''43.39'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '''10.49'
$ws.Range("E13").Value = '  -2.50%  '
$ws.Range("D14").Value = '''4.574.09'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '''3.951.31'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '''14.25'
$ws.Range("E16").Value = '  -3.59%  '
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").Value = '''20.01'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '''69.321.28'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '''437.31'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").Value = '''14.65'
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("D24").Value = '''88.93'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '''12.05'
$ws.Range("E25").Value = '  +5.57%  '
$ws.Range("D26").Value = '''3.84'
$ws.Range("E26").Value = '  +5.11%  '
$ws.Range("D27").Value = '''11.17'
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("D28").Value = '''37.13'
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("D29").Value = '''5.64'
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D31").Value = '''13.43'
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("D32").Value = '''0.129'
$ws.Range("E32").Value = '  -0.71%  '
$ws.Range("E33").Value = '  -1.64%  '
$ws.Range("D34").Value = '''0.455'
$ws.Range("E34").Value = '  +14.61%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '''0.0₃0898'
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '''62.63'
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("D37").Value = '''6.07'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("D38").Value = '''41.14'
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").Value = '''0.996'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("E47").Value = '  +6.33%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '''0.0₆0352'
$ws.Range("E48").Value = '  +2.45%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''2.99'
$ws.Range("E49").Value = '  +4.62%  '
$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").Value = '''3.39'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("E51").Value = '  -2.71%  '
